# Regenerate save_data "K" column (column G) values for Sheet1.
# This mirrors a re-run of the data pipeline that recalculates K (formerly
# "Strike#") per-row and writes the refreshed values back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values keyed by row number (header row is row 1; data rows 2-14).
$newK = @{
    2  = 0
    3  = 1
    4  = 4
    5  = 0
    6  = 2
    7  = 2
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
